$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data rows for "caller_is_employee" (row 6) and "reassignment_count" (row 7)
# so that reassignment_count now appears before caller_is_employee.
for ($c = 1; $c -le 9; $c++) {
    $v6 = $ws.Cells.Item(6, $c).Value2
    $v7 = $ws.Cells.Item(7, $c).Value2
    $ws.Cells.Item(6, $c).Value = $v7
    $ws.Cells.Item(7, $c).Value = $v6
}

# Remove the now-unneeded "predicted_dissatisfaction_delta" column (column I)
$ws.Range("I1:I26").EntireColumn.Delete()
